$wb = $excel.ActiveWorkbook
$wsTraining = $wb.Worksheets.Item(1)
$wsExam = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------------
# 1. Training Dashboard: refresh "PERIOD TO EXPIRE" (H) and "LAST UPDATE" (I)
#    for every training row (3-18) -- the last-update date moved from
#    08-Sep-2025 to 16-Sep-2025, which shrinks the remaining period by 8 days.
# ---------------------------------------------------------------------------
$newPeriods = @{
    3  = 402
    4  = 321
    5  = 324
    6  = 358
    7  = 352
    8  = 359
    9  = 400
    10 = 373
    11 = 377
    12 = 399
    13 = 357
    14 = 388
    15 = 392
    16 = 406
    17 = 405
    18 = 336
}

foreach ($row in 3..18) {
    $wsTraining.Cells.Item($row, 8).Value = $newPeriods[$row]
}

# Write the new "LAST UPDATE" date as literal text (keeping the same
# inline/shared-string cell type the sheet already used) instead of letting
# the date-aware setter coerce it into a serial date. Temporarily force a
# text format, assign the value, then restore the original (General) look by
# copying the format back from the untouched neighbouring STATUS column so
# the style reference collapses back to the shared "body" style.
$dateRange = $wsTraining.Range("I3:I18")
$dateRange.NumberFormat = "@"
$dateRange.Value = "16-Sep-2025"
$wsTraining.Range("J3:J18").Copy()
$dateRange.PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2. Exam Dashboard: wording tweak + wider COMMENTS column so the new text
#    fits.
# ---------------------------------------------------------------------------
$wsExam.Range("E3").Value = "date is valid"
$wsExam.Columns.Item(5).ColumnWidth = 14.1666666666667

# ---------------------------------------------------------------------------
# 3. Header styling: the bold 14pt "title" font and the plain bold "header
#    row" font collapse into a single bold/white font used by both the title
#    banner and the column-header row on each sheet.
# ---------------------------------------------------------------------------
foreach ($ws in @($wsTraining, $wsExam)) {
    $titleFont = $ws.Range("A1").Font
    $titleFont.Size = 11
    $titleFont.Bold = $true
    $titleFont.Color = 16777215

    $headerLastCol = $ws.UsedRange.Columns.Count
    $headerRow = $ws.Range($ws.Cells.Item(2, 1), $ws.Cells.Item(2, $headerLastCol))
    $headerFont = $headerRow.Font
    $headerFont.Bold = $true
    $headerFont.Color = 16777215
}
